$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update phone numbers in B2:B6 to include the "5581" country/area code prefix,
# and format them as plain numbers (General number format, explicitly applied).
$numbers = @(5581945678912, 5581945678912, 5581945678912, 5581945678912, 5581945678912)
for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $numbers[$i]
    $cell.NumberFormat = "General"
}

# Update the selection on the sheet to B2:B6 with active cell B2.
$ws.Range("B2:B6").Select()
$ws.Range("B2").Activate()
